# Timeseries count version 2 - Code Cleaning
# Remove the Karnataka/Bangalore/Mysore rows (old rows 2-4), keeping only
# the Maharashtra data (old rows 5-9), which shifts up to become the new
# rows 2-6. Deleting the rows lets Excel naturally re-flow the remaining
# data and adjust the merged cells (A2:A9 -> A2:A6 style regions) without
# touching any cell styles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:4").Delete()
